$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp footer (row 1)
$ws.Range('A1').Value = 'Datos actualizados a 18 de Junio de 2020 a las 23:03'

# Estados Unidos (row 4): refreshed totals
$ws.Range('B4').Value = 2255289
$ws.Range('C4').Value = 20818
$ws.Range('D4').Value = 923282
$ws.Range('E4').Value = 1211508
$ws.Range('G4').Value = 558
$ws.Range('H4').Value = 120499

# Alemania (row 14): refreshed totals
$ws.Range('B14').Value = 190050
$ws.Range('C14').Value = 546
$ws.Range('E14').Value = 7006
$ws.Range('G14').Value = 17
$ws.Range('H14').Value = 8944

# Egipto overtakes Paises Bajos (rows 29-30 swap places + refreshed totals)
$ws.Range('A29').Value = 'Egipto'
$ws.Range('B29').Value = 50437
$ws.Range('C29').Value = 1218
$ws.Range('D29').Value = 13528
$ws.Range('E29').Value = 34971
$ws.Range('G29').Value = 88
$ws.Range('H29').Value = 1938
$ws.Range('A30').Value = 'Paises Bajos'
$ws.Range('B30').Value = 49319
$ws.Range('C30').Value = 115
$ws.Range('D30').Value = 0
$ws.Range('E30').Value = 0
$ws.Range('G30').Value = 4
$ws.Range('H30').Value = 6078

# Israel (row 50): refreshed totals
$ws.Range('B50').Value = 20036
$ws.Range('C50').Value = 253
$ws.Range('E50').Value = 4215

# Barein (row 51): refreshed totals
$ws.Range('E51').Value = 5721
$ws.Range('G51').Value = 6
$ws.Range('H51').Value = 55

# Costa de Marfil (row 75): refreshed totals
$ws.Range('B75').Value = 6444
$ws.Range('C75').Value = 381
$ws.Range('D75').Value = 2863
$ws.Range('E75').Value = 3532
$ws.Range('G75').Value = 1
$ws.Range('H75').Value = 49

# Mauritania moves ahead of Mayotte/Cuba/Croacia (rows 98-101 shift + refreshed totals)
$ws.Range('A98').Value = 'Mauritania'
$ws.Range('B98').Value = 2424
$ws.Range('C98').Value = 201
$ws.Range('D98').Value = 550
$ws.Range('E98').Value = 1777
$ws.Range('G98').Value = 2
$ws.Range('H98').Value = 97
$ws.Range('A99').Value = 'Mayotte'
$ws.Range('B99').Value = 2383
$ws.Range('C99').Value = 38
$ws.Range('D99').Value = 2066
$ws.Range('E99').Value = 288
$ws.Range('G99').Value = 0
$ws.Range('H99').Value = 29
$ws.Range('A100').Value = 'Cuba'
$ws.Range('B100').Value = 2295
$ws.Range('C100').Value = 15
$ws.Range('D100').Value = 2020
$ws.Range('E100').Value = 190
$ws.Range('G100').Value = 1
$ws.Range('H100').Value = 85
$ws.Range('A101').Value = 'Croacia'
$ws.Range('B101').Value = 2269
$ws.Range('C101').Value = 11
$ws.Range('D101').Value = 2142
$ws.Range('E101').Value = 20
$ws.Range('H101').Value = 107

# Cabo Verde (row 140): refreshed totals
$ws.Range('B140').Value = 688
$ws.Range('C140').Value = 5
$ws.Range('D140').Value = 191
$ws.Range('E140').Value = 485

# Santo Tome y Principe (row 143): refreshed totals
$ws.Range('B143').Value = 646
$ws.Range('C143').Value = 7
$ws.Range('D143').Value = 350
$ws.Range('E143').Value = 294

# Mozambique (row 145): refreshed totals
$ws.Range('B145').Value = 600
$ws.Range('C145').Value = 45
$ws.Range('E145').Value = 182

# Surinam overtakes Birmania (rows 160-161 swap places + refreshed totals)
$ws.Range('A160').Value = 'Surinam'
$ws.Range('B160').Value = 277
$ws.Range('C160').Value = 16
$ws.Range('D160').Value = 74
$ws.Range('E160').Value = 196
$ws.Range('G160').Value = 1
$ws.Range('H160').Value = 7
$ws.Range('A161').Value = 'Birmania'
$ws.Range('B161').Value = 263
$ws.Range('C161').Value = 1
$ws.Range('D161').Value = 187
$ws.Range('E161').Value = 70
$ws.Range('G161').Value = 0
$ws.Range('H161').Value = 6

# Comoras (row 165): refreshed totals
$ws.Range('D165').Value = 141
$ws.Range('E165').Value = 51

# Santa Sede overtakes Islas Turcas y Caicos (rows 208-209 swap places + refreshed totals)
$ws.Range('A208').Value = 'Santa Sede'
$ws.Range('D208').Value = 12
$ws.Range('H208').Value = 0
$ws.Range('A209').Value = 'Islas Turcas y Caicos'
$ws.Range('D209').Value = 11
$ws.Range('H209').Value = 1

# Papua Nueva Guinea overtakes Islas Virgenes Britanicas (rows 213-214 swap places + refreshed totals)
$ws.Range('A213').Value = 'Papua Nueva Guinea'
$ws.Range('D213').Value = 8
$ws.Range('H213').Value = 0
$ws.Range('A214').Value = 'Islas Virgenes Britanicas'
$ws.Range('D214').Value = 7
$ws.Range('H214').Value = 1
